# Apply data updates to worksheet per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E4").Value = 28
$ws.Range("F4").Value = 9
$ws.Range("H4").Value = 23

$ws.Range("E5").Value = 7

$ws.Range("E8").Value = 16

$ws.Range("E12").Value = 8

$ws.Range("E15").Value = 172
$ws.Range("F15").Value = 97
$ws.Range("H15").Value = 138

$ws.Range("E17").Value = 135

$ws.Range("E18").Value = 125

$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 11

$ws.Range("E28").Value = 22

$ws.Range("E29").Value = 19
$ws.Range("F29").Value = 13
$ws.Range("H29").Value = 16

$ws.Range("E37").Value = 60
$ws.Range("F37").Value = 38
$ws.Range("H37").Value = 50

$ws.Range("E38").Value = 84

$ws.Range("E41").Value = 48
$ws.Range("F41").Value = 26
$ws.Range("H41").Value = 37

$ws.Range("E43").Value = 30

$ws.Range("F50").Value = 11
$ws.Range("H50").Value = 19

$ws.Range("E55").Value = 8

$ws.Range("E61").Value = 34

$ws.Range("E62").Value = 50

$ws.Range("E63").Value = 46
$ws.Range("F63").Value = 18
$ws.Range("H63").Value = 26

$ws.Range("E65").Value = 39

$ws.Range("E66").Value = 37
$ws.Range("F66").Value = 26
$ws.Range("H66").Value = 34

$ws.Range("E67").Value = 43
$ws.Range("F67").Value = 25
$ws.Range("H67").Value = 34

$ws.Range("E75").Value = 18
$ws.Range("F75").Value = 9
$ws.Range("H75").Value = 14

$ws.Range("E76").Value = 55
$ws.Range("F76").Value = 21
$ws.Range("H76").Value = 38

$ws.Range("E77").Value = 61
$ws.Range("F77").Value = 24
$ws.Range("H77").Value = 41

$ws.Range("E78").Value = 49

$ws.Range("E80").Value = 32

$ws.Range("E81").Value = 21
$ws.Range("F81").Value = 14
$ws.Range("H81").Value = 19

$ws.Range("E82").Value = 18
$ws.Range("F82").Value = 7
$ws.Range("H82").Value = 13

$ws.Range("E83").Value = 12

$ws.Range("F84").Value = 4
$ws.Range("H84").Value = 5

$ws.Range("E88").Value = 32
$ws.Range("F88").Value = 18
$ws.Range("H88").Value = 26
